$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Regenerate the handoff report: the old GUID-named source file
# (4c672a57-2b8f-4c4a-9f33-64fc334a98f3.md) has been replaced by a new one
# (a47a59d5-9fba-45b8-817a-ccc9dfddf0ed.md), with refreshed xlf content
# hashes and timestamps.
# ---------------------------------------------------------------------------

$oldGuidFile   = "4c672a57-2b8f-4c4a-9f33-64fc334a98f3.md"
$newGuidFile   = "a47a59d5-9fba-45b8-817a-ccc9dfddf0ed.md"
$newDisplay    = "e2e\a47a59d5-9fba-45b8-817a-ccc9dfddf0ed.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# File Name (A2)
$wsOverview.Range("A2").Value = $newGuidFile

# Path And Name (B2) is a hyperlink whose display text must change while
# keeping the same target address / relationship.
$linkB2 = $wsOverview.Range("B2").Hyperlinks.Item(1)
$linkB2.TextToDisplay = $newDisplay
$wsOverview.Range("B2").Hyperlinks.Item(1).Delete()

# Latest HO Xliff Generate Date (G2)
$wsOverview.Range("G2").Value = "2016-08-28 18:56:00"

# --- zh-cn sheet ------------------------------------------------------------
# Source File Name (A2) is a hyperlink whose display text must change while
# keeping the same target address / relationship.
$linkZhA2 = $wsZhCn.Range("A2").Hyperlinks.Item(1)
$linkZhA2.TextToDisplay = $newGuidFile
$wsZhCn.Range("A2").Hyperlinks.Item(1).Delete()

# Latest Handoff File (G2)
$wsZhCn.Range("G2").Value = "a47a59d5-9fba-45b8-817a-ccc9dfddf0ed.d02b17853a42b6786593119d79826fd2e2f85269.zh-cn.xlf"

# Latest Handoff Datetime (H2)
$wsZhCn.Range("H2").Value = "2016-08-28 18:55:56"

# --- de-de sheet ------------------------------------------------------------
# Source File Name (A2) is a hyperlink whose display text must change while
# keeping the same target address / relationship.
$linkDeA2 = $wsDeDe.Range("A2").Hyperlinks.Item(1)
$linkDeA2.TextToDisplay = $newGuidFile
$wsDeDe.Range("A2").Hyperlinks.Item(1).Delete()

# Latest Handoff File (G2)
$wsDeDe.Range("G2").Value = "a47a59d5-9fba-45b8-817a-ccc9dfddf0ed.d02b17853a42b6786593119d79826fd2e2f85269.de-de.xlf"

# Latest Handoff Datetime (H2)
$wsDeDe.Range("H2").Value = "2016-08-28 18:56:00"
